$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match details (columns F:V) between row 22 and row 23 ---
# (Indice/pais/torneio/temporada/data_partida in A:E stay put; only the
# match info that was scraped got reordered.)
$row22 = $ws.Range("F22:V22").Value2
$row23 = $ws.Range("F23:V23").Value2
$ws.Range("F22:V22").Value = $row23
$ws.Range("F23:V23").Value = $row22

# --- Swap match details (columns F:V) between row 30 and row 31 ---
$row30 = $ws.Range("F30:V30").Value2
$row31 = $ws.Range("F31:V31").Value2
$ws.Range("F30:V30").Value = $row31
$ws.Range("F31:V31").Value = $row30

# --- Append new match rows 111-113 (copy formatting from the last
#     existing data row, then overwrite the values) ---

$ws.Range("A110:V110").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)
$ws.Range("A110:V110").Copy()
$ws.Range("A112:V112").PasteSpecial(-4122)
$ws.Range("A110:V110").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)

# Row 111
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = "belgium"
$ws.Cells.Item(111, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(111, 4).Value = "2023-2024"
$ws.Cells.Item(111, 5).Value = 45242.66666666666
$ws.Cells.Item(111, 6).Value = "Genk"
$ws.Cells.Item(111, 7).Value = 3
$ws.Cells.Item(111, 8).Value = "Leuven"
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 1.36
$ws.Cells.Item(111, 11).Value = "05/11/2023 16:12"
$ws.Cells.Item(111, 12).Value = 1.33
$ws.Cells.Item(111, 13).Value = "12/11/2023 15:36"
$ws.Cells.Item(111, 14).Value = 5.78
$ws.Cells.Item(111, 15).Value = "05/11/2023 16:12"
$ws.Cells.Item(111, 16).Value = 6.09
$ws.Cells.Item(111, 17).Value = "12/11/2023 15:36"
$ws.Cells.Item(111, 18).Value = 7.04
$ws.Cells.Item(111, 19).Value = "05/11/2023 16:12"
$ws.Cells.Item(111, 20).Value = 8.140000000000001
$ws.Cells.Item(111, 21).Value = "12/11/2023 15:36"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-leuven/pAqHIgXK/"

# Row 112
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "belgium"
$ws.Cells.Item(112, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(112, 4).Value = "2023-2024"
$ws.Cells.Item(112, 5).Value = 45242.77083333334
$ws.Cells.Item(112, 6).Value = "Gent"
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = "Anderlecht"
$ws.Cells.Item(112, 9).Value = 1
$ws.Cells.Item(112, 10).Value = 1.83
$ws.Cells.Item(112, 11).Value = "05/11/2023 19:43"
$ws.Cells.Item(112, 12).Value = 2.03
$ws.Cells.Item(112, 13).Value = "12/11/2023 18:05"
$ws.Cells.Item(112, 14).Value = 3.95
$ws.Cells.Item(112, 15).Value = "05/11/2023 19:43"
$ws.Cells.Item(112, 16).Value = 3.63
$ws.Cells.Item(112, 17).Value = "12/11/2023 18:29"
$ws.Cells.Item(112, 18).Value = 4.04
$ws.Cells.Item(112, 19).Value = "05/11/2023 19:43"
$ws.Cells.Item(112, 20).Value = 3.72
$ws.Cells.Item(112, 21).Value = "12/11/2023 18:29"
$ws.Cells.Item(112, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-anderlecht/4fQmDXfr/"

# Row 113
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = "belgium"
$ws.Cells.Item(113, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(113, 4).Value = "2023-2024"
$ws.Cells.Item(113, 5).Value = 45242.80208333334
$ws.Cells.Item(113, 6).Value = "Royale Union SG"
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = "Kortrijk"
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1.23
$ws.Cells.Item(113, 11).Value = "05/11/2023 19:43"
$ws.Cells.Item(113, 12).Value = 1.19
$ws.Cells.Item(113, 13).Value = "12/11/2023 19:14"
$ws.Cells.Item(113, 14).Value = 6.6
$ws.Cells.Item(113, 15).Value = "05/11/2023 19:43"
$ws.Cells.Item(113, 16).Value = 7.51
$ws.Cells.Item(113, 17).Value = "12/11/2023 19:14"
$ws.Cells.Item(113, 18).Value = 9.119999999999999
$ws.Cells.Item(113, 19).Value = "05/11/2023 19:43"
$ws.Cells.Item(113, 20).Value = 14.26
$ws.Cells.Item(113, 21).Value = "12/11/2023 19:14"
$ws.Cells.Item(113, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/royale-union-sg-kortrijk/jmpDJZHE/"
